$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "329.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.34%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.04%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.579"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.57%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08063"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.10%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.983"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.55%"

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.574"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-4.58%"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9519"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.56%"

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1159"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.35%"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1858"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.73%"

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "11.86"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "38.77%"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09879"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.07%"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04747"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "14.27%"

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1069"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.26%"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001287"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.51%"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04239"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.38%"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005938"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.34%"

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.372"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-5.17%"

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.328"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.92%"

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3475"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.29%"

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1412"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.38%"

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2507"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.27%"

$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001251"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.38%"

$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004338"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.15%"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.35%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.49%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02645"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "0.02%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05554"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.65%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007572"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.04%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1407"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.39%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.008081"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-28.05%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.14%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008867"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-8.29%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007090"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.20%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.16%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.002302"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "1.17%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.004852"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "35.94%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.16%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.16%"
